# Refitting NCDEs to individual patients (for manuscript figure)
# Adds a "Label" column (H) marking Control (0) vs MDD (1) rows, and
# updates the refitted D/E/F prediction values produced by the new fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Label" header in H1 (match formatting of the other headers) ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Label"

# --- Updated refitted values (block 1: Iterations = 100, rows 2-11) ---
$ws.Range("D2").Value = 0.6641409288527483
$ws.Range("E2").Value = 0.6641409288527483

$ws.Range("D5").Value = 0.5531318950330061
$ws.Range("E5").Value = 0.5531318950330061

$ws.Range("D6").Value = 0.1654520705912127
$ws.Range("E6").Value = 0.1654520705912127

$ws.Range("D7").Value = 0.8507932340944957
$ws.Range("E7").Value = 0.1492067659055043

$ws.Range("D8").Value = 0.6590116282705859
$ws.Range("E8").Value = 0.3409883717294141

$ws.Range("D10").Value = 0.6536332890459873
$ws.Range("E10").Value = 0.3463667109540127

$ws.Range("F11").Value = 0.6166301965713501

# --- New "Label" column values ---
# Block 1 (rows 2-11): Controls = 0, MDD = 1
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 1
$ws.Range("H8").Value = 1
$ws.Range("H9").Value = 1
$ws.Range("H10").Value = 1
$ws.Range("H11").Value = 1

# Block 2 (rows 12-21): Controls = 0, MDD = 1
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
